$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MASTER")

# Update EXECUTE column values for TC_1019 / TC_1020 rows from "Yes" to "No"
$ws.Range("C20").Value = "No"
$ws.Range("C21").Value = "No"

# Update the active selection to D21
$ws.Range("D21").Select()
